$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Sheet3: update the rolling lookup table (A20:B36) with the figures
#    recomputed after adding the "08-nov" day of data. Rows 21/22/25
#    (CHEETOS 94GRX24, DORITOS QUESO 140GX19, LAYS CEBOLLA CARAMELIZADA)
#    are unaffected and stay at 0.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Range("B20").Value = 12.072133765285232
$ws3.Range("B23").Value = 6.0550147849501252
$ws3.Range("B24").Value = 5.9556439858799726
$ws3.Range("B26").Value = 5.2883303828142978
$ws3.Range("B27").Value = 14.674201013192437
$ws3.Range("B28").Value = 3.7711403702315951
$ws3.Range("B29").Value = 4.0718899891423241
$ws3.Range("B30").Value = 3.9680079980003136
$ws3.Range("B31").Value = 13.050987549754897
$ws3.Range("B32").Value = 3.9013361355079228
$ws3.Range("B33").Value = 16.445707798319088
$ws3.Range("B34").Value = 8.0396707445872959
$ws3.Range("B35").Value = 4.9376131140857851
$ws3.Range("B36").Value = 47.942381785227276

# (Sheet3 C2:C18 VLOOKUP formulas, and Sheet1 CB/CC VLOOKUP formulas that
#  reference them, recalculate automatically.)

# ---------------------------------------------------------------------
# 2. Sheet1: append the new "08-nov" day column (CK), copying the
#    now-current VLOOKUP figures the same way the previous day's column
#    (CJ) had frozen the prior day's figures as plain values.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("CK1").Value = "08-nov"

$ws1.Range("CK2").Value = 3.7711403702315951
$ws1.Range("CK3").Value = 4.0718899891423241
$ws1.Range("CK4").Value = 5.2883303828142978
$ws1.Range("CK5").Value = 14.674201013192437
$ws1.Range("CK6").Value = 6.0550147849501252
$ws1.Range("CK7").Value = 5.9556439858799726
$ws1.Range("CK8").Value = 0
$ws1.Range("CK9").Value = 4.9376131140857851
$ws1.Range("CK10").Value = 16.445707798319088
$ws1.Range("CK11").Value = 8.0396707445872959
$ws1.Range("CK12").Value = 12.072133765285232
$ws1.Range("CK13").Value = 0
$ws1.Range("CK14").Value = 47.942381785227276
$ws1.Range("CK15").Value = 3.9013361355079228
$ws1.Range("CK16").Value = 0
$ws1.Range("CK17").Value = 3.9680079980003136
$ws1.Range("CK18").Value = 13.050987549754897

# Match the number formatting used by the rest of the day columns
# (integer display format, same as CJ2:CJ18).
$ws1.Range("CK2:CK18").NumberFormat = "0"

# ---------------------------------------------------------------------
# 3. Housekeeping that Excel performs when the new column pushes the
#    previous helper columns out of view: CA:CB lose their "best fit"
#    auto width and get hidden along with the newly-unused CC:CF block.
# ---------------------------------------------------------------------
$ws1.Range("CA1:CB1").EntireColumn.Hidden = $true
$ws1.Range("CC1:CF1").EntireColumn.ColumnWidth = 0
$ws1.Range("CC1:CF1").EntireColumn.Hidden = $true

# Leave the selection on the newly entered cell, as the author did.
$ws1.Range("CK2").Select()
